# Updates crypto Price (column D) and Volume(1h) (column E) cells
# with refreshed data, per the Sun Feb 11 09:25:12 UTC 2024 GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.225.88"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.526.90"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.93"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.50"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.14"
$ws.Range("E10").Value = "  +5.86%  "
$ws.Range("E11").Value = "  +11.62%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "2.922.03"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "2.531.21"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "48.061.88"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.24"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.46"
$ws.Range("E24").Value = "  +9.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.27"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.98"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.65"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.03"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0796"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.28"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "120.03"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0302"
$ws.Range("D45").Value = "2.023.65"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  +4.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").Value = "  +7.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.28"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.82"
$ws.Range("E51").Value = "  +2.92%  "
